# Update notebook/provsvar_tabel.xlsx datasets (rows 2-4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = 8201113194
$ws.Range("C2").Value = "30/11/2021"
$ws.Range("E2").Value = 18
$ws.Range("H2").Value = 1
$ws.Range("K2").Value = 0.3
$ws.Range("N2").Value = 0.39
$ws.Range("Q2").Value = 0.34
$ws.Range("S2").Value = "RÖD"
$ws.Range("S2").Interior.Color = 255

# --- Row 3 ---
$ws.Range("A3").Value = 23
$ws.Range("B3").Value = 8201113197
$ws.Range("C3").Value = "29/11/2021"
$ws.Range("E3").Value = 16
$ws.Range("G3").Value = "T3, P-"
$ws.Range("H3").Value = 0.87
$ws.Range("K3").Value = 6.4
$ws.Range("L3").Value = "GUL"
$ws.Range("L3").Interior.Color = 65535
$ws.Range("N3").Value = 0.52
$ws.Range("Q3").Value = 0.48
$ws.Range("R3").Value = "GRON"
$ws.Range("R3").Interior.Color = 32768

# --- Row 4 ---
$ws.Range("A4").Value = 26
$ws.Range("B4").Value = 8201113206
$ws.Range("C4").Value = "26/11/2021"
$ws.Range("E4").Value = 16
$ws.Range("G4").Value = "T3, P-"
$ws.Range("H4").Value = 1.2
$ws.Range("N4").Value = 0.4
$ws.Range("Q4").Value = 0.18
